$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TwitterUserLoginTest")

# Remove the hyperlink on A2 (mailto link tied to the username value)
$ws.Hyperlinks.Delete()

# Clear the username/password test data (and the now-unused shared strings
# go away automatically when the workbook is saved)
$ws.Range("A2:B2").ClearContents()

# Make this sheet the active tab, with D7 selected
$ws.Activate()
$ws.Range("D7").Select()
